$p = $ppt.ActivePresentation

# --- Slide 16 ("DEMO") subtitle shortened ---
# Before: "Managing Users with the Azure ARM REST API"
#  After: "Managing Users"
# (Keeps the existing "Managing " run intact and replaces the remainder
#  of the text with "Users" so the trailing run picks up the formatting
#  of the old "API" run, matching a normal type-over-selection edit.)
$slide = $p.Slides.Item(16)
$shape = $slide.Shapes("Text Placeholder 4")
$tr = $shape.TextFrame.TextRange

$prefix = "Managing "
if ($tr.Text.StartsWith($prefix)) {
    $tailStart = $prefix.Length + 1
    $tailLen = $tr.Length - $prefix.Length
    $tail = $tr.Characters($tailStart, $tailLen)
    $tail.Text = "Users"
}
